# Update the list of defined vulnerabilities:
#  - insert a new "MLDNVERDEV2" row right after the existing MLDNVERDEV row
#  - fix up the "ID" column on the row that ends up at r36 (ICMP broadcast-echo
#    network row), which should read "Node" rather than "Network"
#  - append two new rows describing the new ICMPv6 Redirect vulnerability

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10 (pushes the existing rows 10-36 down to 11-37,
# carrying their formatting - e.g. the wrap-text style on some Code cells -
# along with them).
$ws.Rows("10:10").Insert()

$ws.Cells.Item(10, 1).Value = "Node"
$ws.Cells.Item(10, 2).Value = "a,a+"
$ws.Cells.Item(10, 3).Value = "PTV-NET-IDENT-ACTIVE-MLDNVERDEV2"
$ws.Cells.Item(10, 4).Value = "Device only responds to illegitimate MLDv1 queries even though MLDv2 queries are sent, possibly downgraded"

# After the shift, the row that lands at 36 (PTV-NET-IDENT-ICMP-BRCASTECHO) needs
# its ID column corrected to "Node".
$ws.Cells.Item(36, 1).Value = "Node"

# Append the two new rows describing the ICMPv6 Redirect vulnerability. Fill
# column-by-column (both codes, then both descriptions) to mirror how the
# shared-string table ended up ordered in the authored workbook.
$ws.Cells.Item(38, 1).Value = "Network"
$ws.Cells.Item(39, 1).Value = "Node"

$ws.Cells.Item(38, 2).Value = "p,a,a+"
$ws.Cells.Item(39, 2).Value = "p,a,a+"

$ws.Cells.Item(38, 3).Value = "PTV-NET-MITM-ICMP6REDIR"
$ws.Cells.Item(39, 3).Value = "PTV-NET-MITM-ICMP6REDIRDEV"

$ws.Cells.Item(38, 4).Value = "Network does not block ICMPv6 Redirect messages"
$ws.Cells.Item(39, 4).Value = "Device communication can be redirected using ICMPv6 Redirect"

# Match the final selection/view state from the edit.
$ws.Range("K36").Select()
